$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Property1")
$ws2 = $wb.Worksheets.Item("Record_Cooldown")

# --- Property1!G3: was TRUE(1) in the "Save" row, now FALSE(0); also its
#     cell style reverts to the "no-bottom-border" look shared by rows 7-9
#     (copy format from G7, which already carries that style, then set value).
$ws1.Range("G7").Copy()
$ws1.Range("G3").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("G3").Value = 0

# --- Property1!B6:G6 ("Cache" row): every flag flips to TRUE(1), and the
#     whole row is restyled to match the "full border" look used by rows 2-4
#     (copy format from the already-correctly-styled C6, then set values).
$ws1.Range("C6").Copy()
$ws1.Range("B6").PasteSpecial(-4122)
$ws1.Range("G6").PasteSpecial(-4122)

$ws1.Range("B6").Value = 1
$ws1.Range("C6").Value = 1
$ws1.Range("D6").Value = 1
$ws1.Range("E6").Value = 1
$ws1.Range("F6").Value = 1
$ws1.Range("G6").Value = 1

# --- Extend the TRUE/FALSE list validation to the newly-editable B6:E6 cells.
$ws1.Range("B6:E6").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# --- Selection / active-sheet bookkeeping: the user ends up on Property1
#     with G3 selected (instead of Record_Cooldown being the active tab).
[void]$ws1.Range("G3").Select()
